$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mã học sinh" (student code) column moves from F to D (right after
# "Số điện thoại"), pushing "Ngày sinh" (was D) to E and "Giới tính"
# (was E) to F. Column F's old content (the old student-code sample data)
# is discarded - it is replaced by new sample values in D.
#
# Work right-to-left so each shift reads the not-yet-overwritten source
# column (format, then values, copied separately so blank source cells
# do not spuriously create a styled destination cell where none should
# exist).

# Rows 1-5: F (Giới tính) <- old E, then E (Ngày sinh) <- old D
$ws.Range("E1:E5").Copy()
$ws.Range("F1:F5").PasteSpecial(-4122)
$ws.Range("E1:E5").Copy()
$ws.Range("F1:F5").PasteSpecial(-4163)

$ws.Range("D1:D5").Copy()
$ws.Range("E1:E5").PasteSpecial(-4122)
$ws.Range("D1:D5").Copy()
$ws.Range("E1:E5").PasteSpecial(-4163)

# Row 6 only ever had an (empty) D6 - it becomes an empty E6, no F6.
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("D6").ClearContents()

# D takes on the number-format/fill style that column C (Số điện thoại)
# uses, since the new "Mã học sinh" column sits right next to it.
$ws.Range("C1:C6").Copy()
$ws.Range("D1:D6").PasteSpecial(-4122)

$ws.Application.CutCopyMode() = $false

# New header + sample data for the relocated "Mã học sinh" column.
$ws.Range("D1").Value() = "Mã học sinh"
$ws.Range("D2").Value() = "001"
$ws.Range("D3").Value() = "002"
$ws.Range("D4").Value() = "003"
$ws.Range("D5").Value() = "004"
$ws.Range("D6").Value() = ""

# The helper note in row 3 (merged G3:N3) is rewritten to the shorter
# "can be left blank" wording.
$ws.Range("G3").Value() = "*Mã học sinh: có thể để trống"

# Selection moved to D12.
$ws.Range("D12").Select()
